$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, pushing the existing data (rows 7-56) down by one
# row (to rows 8-57). Excel auto-updates the sheet dimension (A1:R56 -> A1:R57).
$ws.Rows.Item(7).Insert()

# The row that is now at index 8 holds what used to be row 7's data. Copy that
# row's values into the newly-inserted (blank) row 7 so the new record starts
# out identical to the record it is based on.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(7, $col).Value = $ws.Cells.Item(8, $col).Value()
}

# Now set the new row's date (column D) to the new reporting date.
$newDate = Get-Date -Year 2021 -Month 10 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(7, 4).Value = $newDate
